$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44559
$ws.Range("K2").Value = 'Modesto'
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25500
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1417
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44559
$ws.Range("K3").Value = 'Modesto'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 320
$ws.Range("N3").Value = 22000
$ws.Range("O3").Value = 23000
$ws.Range("P3").Value = 22500
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1250
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44580
$ws.Range("K4").Value = 'Modesto'
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 22500
$ws.Range("O4").Value = 23000
$ws.Range("P4").Value = 22750
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 1264
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44580
$ws.Range("K5").Value = 'Modesto'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 19500
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19750
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1097
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44175
$ws.Range("K6").Value = 'Castle Brite'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("Q6").Value = '$/caja 18 kilos'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 1194
$ws.Range("T6").Value = 18

# Row 7
$ws.Range("D7").Value = 44545
$ws.Range("K7").Value = 'Castle Brite'
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 340
$ws.Range("N7").Value = 22500
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22750
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1264
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44545
$ws.Range("K8").Value = 'Castle Brite'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 400
$ws.Range("N8").Value = 20500
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20750
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1153
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44545
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 15500
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15750
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 875
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44566
$ws.Range("K10").Value = 'Modesto'
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1306
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44566
$ws.Range("K11").Value = 'Modesto'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 160
$ws.Range("N11").Value = 21000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 21500
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1194
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44552
$ws.Range("K12").Value = 'Castle Brite'
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 360
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 1139
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44552
$ws.Range("K13").Value = 'Castle Brite'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 280
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 19000
$ws.Range("P13").Value = 18500
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 1028
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 44161
$ws.Range("K14").Value = 'Dina'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20500
$ws.Range("P14").Value = 20250
$ws.Range("Q14").Value = '$/caja 15 kilos'
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 1350
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("D15").Value = 44161
$ws.Range("K15").Value = 'Dina'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 18500
$ws.Range("P15").Value = 18250
$ws.Range("Q15").Value = '$/caja 15 kilos'
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("S15").Value = 1217
$ws.Range("T15").Value = 15

# Row 16
$ws.Range("D16").Value = 44573
$ws.Range("K16").Value = 'Modesto'
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 20500
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20750
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 1153
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44573
$ws.Range("K17").Value = 'Modesto'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = 17500
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17750
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 986
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44160
$ws.Range("K18").Value = 'Castle Brite'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 240
$ws.Range("N18").Value = 20500
$ws.Range("O18").Value = 21000
$ws.Range("P18").Value = 20750
$ws.Range("Q18").Value = '$/caja 15 kilos'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 1383
$ws.Range("T18").Value = 15

# Row 19
$ws.Range("D19").Value = 44553
$ws.Range("K19").Value = 'Modesto'
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 360
$ws.Range("N19").Value = 23000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 23500
$ws.Range("Q19").Value = '$/caja 16 kilos'
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 1469
$ws.Range("T19").Value = 16

# Row 20
$ws.Range("D20").Value = 44553
$ws.Range("K20").Value = 'Modesto'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = 21000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 21500
$ws.Range("Q20").Value = '$/caja 16 kilos'
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 1344
$ws.Range("T20").Value = 16

# Row 21
$ws.Range("D21").Value = 44553
$ws.Range("K21").Value = 'Modesto'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 240
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 17500
$ws.Range("Q21").Value = '$/caja 16 kilos'
$ws.Range("R21").Value = 'Región Metropolitana'
$ws.Range("S21").Value = 1094
$ws.Range("T21").Value = 16

# Row 22
$ws.Range("D22").Value = 44189
$ws.Range("K22").Value = 'Dina'
$ws.Range("L22").Value = 'Especial'
$ws.Range("M22").Value = 120
$ws.Range("N22").Value = 23500
$ws.Range("O22").Value = 24000
$ws.Range("P22").Value = 23750
$ws.Range("Q22").Value = '$/caja 18 kilos'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1319
$ws.Range("T22").Value = 18

# Row 23
$ws.Range("D23").Value = 44189
$ws.Range("K23").Value = 'Dina'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 21500
$ws.Range("O23").Value = 22000
$ws.Range("P23").Value = 21750
$ws.Range("Q23").Value = '$/caja 18 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 1208
$ws.Range("T23").Value = 18

# Row 24
$ws.Range("D24").Value = 44546
$ws.Range("K24").Value = 'Castle Brite'
$ws.Range("L24").Value = 'Especial'
$ws.Range("M24").Value = 300
$ws.Range("N24").Value = 22500
$ws.Range("O24").Value = 23000
$ws.Range("P24").Value = 22750
$ws.Range("Q24").Value = '$/caja 18 kilos'
$ws.Range("R24").Value = 'Región Metropolitana'
$ws.Range("S24").Value = 1264
$ws.Range("T24").Value = 18

# Row 25
$ws.Range("D25").Value = 44546
$ws.Range("K25").Value = 'Castle Brite'
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 20500
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 20750
$ws.Range("Q25").Value = '$/caja 18 kilos'
$ws.Range("R25").Value = 'Región Metropolitana'
$ws.Range("S25").Value = 1153
$ws.Range("T25").Value = 18
